{"js": "// Remove the leftover \"vnpt.SiteAddress\" placeholder run that follows the\n// \"\u0110\u1ecba ch\u1ec9: \" label in the \"B\u00ean A\" block (M\u1eabu 26 template clean-up).\nconst results = context.document.body.search(\"vnpt.SiteAddress\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the leftover \"vnpt.SiteAddress\" placeholder run that follows the\n# \"\u0110\u1ecba ch\u1ec9: \" label in the \"B\u00ean A\" block (M\u1eabu 26 template clean-up).\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$range.Find.ClearFormatting()\n$range.Find.Replacement.ClearFormatting()\n$range.Find.Text = \"vnpt.SiteAddress\"\n$range.Find.Replacement.Text = \"\"\n$range.Find.Forward = $true\n$range.Find.Wrap = 1  # wdFindContinue\n\n$range.Find.Execute($range.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $range.Find.Replacement.Text, 2) | Out-Null\n"}
